$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new data rows appended after the existing last row (33), extending
# the sheet to A1:H36.
$newRows = @(
    @("Address Flow 4", "City Flow 4", " ", "Flow 4 First", "Flow 4 Last", " ", "444", "4444"),
    @("Flow 6 Addresst", "Flow 6 City", " ", "Flow 6 First Name", "Flow 7 Lasst Name", " ", "666", "6666"),
    @("Address Flow 1", "City Flow 1", " ", "first flow 1", "last flow 2", " ", "111", "1111")
)

# Values that are purely digits (e.g. "444") must still be stored as TEXT,
# matching every other cell in this column throughout the sheet. A plain
# `.Value = "444"` assignment is auto-coerced to a number by Excel, and
# forcing text via NumberFormat="@" (or a leading apostrophe) leaves a
# permanent "quote prefix"/text-format style on the cell that the source
# file never had. Instead, enter it as a formula that evaluates to a text
# string, then flatten it to a plain value via copy/paste-special (values
# only) -- this keeps the text type without touching the cell's style.
function Set-TextValue($cell, [string]$text) {
    if ($text -match '^-?\d+(\.\d+)?$') {
        $escaped = $text.Replace('"', '""')
        $cell.Formula = '="' + $escaped + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    } else {
        $cell.Value = $text
    }
}

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $newRows[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $cell = $ws.Cells.Item($rowNum, $col)
        Set-TextValue $cell $rowData[$col - 1]
    }
}
